$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Arkusz1" (sheet1): update headcounts / assignees (risk mitigation:
# KaJa, PiSm -> swap people out) and a new L1 note.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("B1").Value = "6+ (4 os.)"
$ws1.Range("D1").Value = "8+ (4 os.)"
$ws1.Range("F1").Value = "10+ (4 os.)"
$ws1.Range("L1").Value = "Piotr S. ??/ Lukasz L."

$ws1.Range("C2").Value = "Tomek,`r`nDawid."
$ws1.Range("E2").Value = "Piotrek S. ??`r`nPiorek K.`r`nMarcin"
$ws1.Range("G2").Value = "Krzysiu, `r`nZuza"
$ws1.Range("I2").Value = "Lukasz G.`r`nJacek"

$ws1.Range("C3").Value = "Tomek,`r`nLukasz G. "
$ws1.Range("G3").Value = "Piotrek K,`r`nDawid,`r`nLukasz L. ??"

$ws1.Columns.Item(12).ColumnWidth = 23.7109375

# ---------------------------------------------------------------------------
# Sheet "Arkusz2" (sheet2): build a brand-new "Co / Kto / Uwagi" task table.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Columns.Item(1).ColumnWidth = 55.7109375
$ws2.Columns.Item(2).ColumnWidth = 19
$ws2.Columns.Item(3).ColumnWidth = 32.140625

$ws2.Range("A1").Value = "Co"
$ws2.Range("B1").Value = "Kto"
$ws2.Range("C1").Value = "Uwagi"
$ws2.Range("A1:C1").Style = "Good"

$ws2.Range("A2").Value = "zakupy pod lutowanie:"
$ws2.Range("B2").Value = "Marcin M."

$ws2.Range("A3").Value = "zakupy catering:"
$ws2.Range("B3").Value = "Lukasz L."
$ws2.Range("C3").Value = "'+ picie do obiadu osobno"

$ws2.Range("A4").Value = "sphero:"
$ws2.Range("B4").Value = "Piotrek K"

$ws2.Range("A5").Value = "robogames:"
$ws2.Range("B5").Value = "Tomek"

$ws2.Range("A6").Value = "gadżety/upominki"
$ws2.Range("B6").Value = "Piotrek S./Zuza"

$ws2.Range("A7").Value = "przygotowanie sal (piątek):"
$ws2.Range("B7").Value = "wszyscy"

$ws2.Range("A8").Value = "kompy:"
$ws2.Range("B8").Value = "wszyscy"

$ws2.Range("A9").Value = "soft:"
$ws2.Range("B9").Value = "Marcin M."

$ws2.Range("A10").Value = "iPady:"
$ws2.Range("B10").Value = "community (1 szt.), Piotrek K"

$ws2.Range("A11").Value = "tablety:"
$ws2.Range("B11").Value = "community (1 szt.) + Lukasz L."

$ws2.Range("A12").Value = "LEGO:"
$ws2.Range("B12").Value = "Dawid + LuLa (demo)"

$ws2.Range("A13").Value = "arduino:"
$ws2.Range("B13").Value = "Jacek"

$ws2.Range("A14").Value = "zadnia na pendrive'y:"
$ws2.Range("B14").Value = "Krzysiu"

$ws2.Range("A15").Value = "etykietki:"
$ws2.Range("B15").Value = "Lukasz L."

$ws2.Range("A16").Value = "sale bookowanie (PLUMy -> Piotr S.)/rzutniki/podzial sal"
$ws2.Range("B16").Value = "Lukasz L."
$ws2.Range("C16").Value = "PLUM 1+2, PLUM3 (zglosic Dorocie)"

$ws2.Range("A17").Value = "dyplomy:"
$ws2.Range("B17").Value = "Piotr S."

$ws2.Range("A18").Value = "baterie / ladowarka"
$ws2.Range("B18").Value = "Lukasz L."

$ws2.Range("A19").Value = "zgloszenie event'u"
$ws2.Range("B19").Value = "Lukasz L."
$ws2.Range("C19").Value = "Ochrona/Magda"

$ws2.Range("A20").Value = "pizza: "
$ws2.Range("B20").Value = "Koordynator "
$ws2.Range("C20").Value = "proste menu "

# Remember sheet2's own selection (C16) without stealing the active tab,
# then re-select sheet1's L1 last so Arkusz1 stays the active/visible tab.
$ws2.Range("C16").Select()
$ws1.Range("L1").Select()
